$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.167686939239502
$ws.Range("B1").Value = 2.340146064758301
$ws.Range("D1").Value = 1.449784636497498
$ws.Range("E1").Value = 0.9424439072608948
